# Apply the edits described in the commit:
#  1. Update the text of the shared string used by A1 (append warning).
#  2. Reset the sheet's selection back to A1 (the default), clearing the
#     previously-saved "A2" selection that was stored in the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update cell text
$ws.Range("A1").Value = "Delete this file and add your files here. Do it asap."

# 2) Reset the active selection to A1 (removing the stale A2 selection)
$ws.Activate()
$ws.Range("A1").Select()
